$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (losing trailing zeros / exact formatting),
# matching the original inlineStr (text) storage used in the workbook.
$textCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D14", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D27", "D29", "D30", "D31", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D41", "D42", "D43", "D44", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated coin data (price & 1h volume/change columns, plus the few
# rows whose coin name/link/price/volume were swapped or changed).
$ws.Range("D2").Value = "56.452.30"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "3.003.00"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.33%  "
$ws.Range("D5").Value = "506.31"
$ws.Range("E5").Value = "  +1.31%  "
$ws.Range("D6").Value = "138.86"
$ws.Range("E6").Value = "  +0.72%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").Value = "0.432"
$ws.Range("E8").Value = "  +0.87%  "
$ws.Range("D9").Value = "7.14"
$ws.Range("E9").Value = "  -2.60%  "
$ws.Range("E11").Value = "  +2.89%  "
$ws.Range("D12").Value = "3.474.81"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("E13").Value = "  -1.94%  "
$ws.Range("D14").Value = "26.30"
$ws.Range("E14").Value = "  +1.73%  "
$ws.Range("E15").Value = "  +2.19%  "
$ws.Range("D16").Value = "56.294.63"
$ws.Range("E16").Value = "  -1.14%  "
$ws.Range("D17").Value = "6.05"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "2.970.17"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "13.07"
$ws.Range("E19").Value = "  +3.76%  "
$ws.Range("D20").Value = "8.02"
$ws.Range("E20").Value = "  +2.61%  "
$ws.Range("D21").Value = "329.36"
$ws.Range("E21").Value = "  +3.29%  "
$ws.Range("D22").Value = "0.992"
$ws.Range("E22").Value = "  -0.65%  "
$ws.Range("D23").Value = "0.497"
$ws.Range("E23").Value = "  +2.56%  "
$ws.Range("D24").Value = "64.48"
$ws.Range("E24").Value = "  +2.14%  "
$ws.Range("D25").Value = "3.094.75"
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("D27").Value = "0.163"
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("D28").Value = "0.0₃0908"
$ws.Range("E28").Value = "  +1.95%  "
$ws.Range("D29").Value = "6.54"
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("D30").Value = "7.05"
$ws.Range("E30").Value = "  -0.21%  "
$ws.Range("D31").Value = "1.79"
$ws.Range("E31").Value = "  +1.72%  "
$ws.Range("E32").Value = "  +1.04%  "
$ws.Range("D33").Value = "20.34"
$ws.Range("E33").Value = "  +1.20%  "
$ws.Range("D34").Value = "151.96"
$ws.Range("E34").Value = "  -1.77%  "
$ws.Range("D35").Value = "4.58"
$ws.Range("E35").Value = "  -0.48%  "
$ws.Range("D36").Value = "5.80"
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("D37").Value = "25.50"
$ws.Range("D38").Value = "1.27"
$ws.Range("E38").Value = "  +1.75%  "
$ws.Range("D39").Value = "0.0664"
$ws.Range("E39").Value = "  -0.22%  "
$ws.Range("D40").Value = "2.997.99"
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("D41").Value = "36.71"
$ws.Range("E41").Value = "  -2.44%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "3.82"
$ws.Range("E43").Value = "  +2.57%  "
$ws.Range("D44").Value = "0.656"
$ws.Range("E44").Value = "  +2.77%  "
$ws.Range("D45").Value = "2.183.43"
$ws.Range("E45").Value = "  -0.84%  "
$ws.Range("E46").Value = "  -1.40%  "
$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D47").Value = "5.95"
$ws.Range("E47").Value = "  +0.32%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "0.939"
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("D49").Value = "19.86"
$ws.Range("E49").Value = "  +3.39%  "
$ws.Range("D50").Value = "0.0238"
$ws.Range("E50").Value = "  +1.51%  "
$ws.Range("D51").Value = "0.0859"
$ws.Range("E51").Value = "  -1.98%  "

# Restore default (unstyled) cell style on the cells where we forced text
# formatting above, so formatting matches the original (unstyled) cells.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
